$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "and"
$ws.Range("A5").Value = "smart"
$ws.Range("A6").Value = "and stuff"

$ws.Range("B5").Select()
